# Auto-generated Excel COM-interop script to apply Golem_Profits data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2555.5557
$ws.Range("I17").Value = 2833.3333
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 8499.999899999999
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -8331.999899999999
$ws.Range("N17").Value = -6336
$ws.Range("H26").Value = 2000
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 2000
$ws.Range("M26").Value = -1656
$ws.Range("H34").Value = 4999
$ws.Range("J34").Value = 4999
$ws.Range("L34").Value = 4999
$ws.Range("N34").Value = -5405
$ws.Range("H36").Value = 4999
$ws.Range("J36").Value = 4999
$ws.Range("L36").Value = 4999
$ws.Range("N36").Value = -6429
$ws.Range("H47").Value = 40790
$ws.Range("J47").Value = 50000
$ws.Range("L47").Value = 50000
$ws.Range("N47").Value = -51944
$ws.Range("H54").Value = 1076
$ws.Range("I54").Value = 1076
$ws.Range("K54").Value = 1076
$ws.Range("M54").Value = -590
$ws.Range("H70").Value = 2845.24
$ws.Range("I70").Value = 2845.24
$ws.Range("K70").Value = 8535.719999999999
$ws.Range("M70").Value = -8265.719999999999
$ws.Range("H73").Value = 2845.24
$ws.Range("I73").Value = 2845.24
$ws.Range("K73").Value = 8535.719999999999
$ws.Range("M73").Value = -7599.719999999999
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H124").Value = 17500
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 12684.333
$ws.Range("J39").Value = 12684.333
$ws.Range("L39").Value = 12684.333
$ws.Range("N39").Value = -13462.333
$ws.Range("H94").Value = 2369.1538
$ws.Range("I94").Value = 1974.875
$ws.Range("K94").Value = 1974.875
$ws.Range("M94").Value = -1523.875
$ws.Range("H99").Value = 2854.5
$ws.Range("I99").Value = 2714
$ws.Range("J99").Value = 4400
$ws.Range("K99").Value = 2714
$ws.Range("L99").Value = 4400
$ws.Range("M99").Value = -1216
$ws.Range("N99").Value = -7396

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1998.5714
$ws.Range("I32").Value = 1998.5714
$ws.Range("K32").Value = 1998.5714
$ws.Range("M32").Value = -1682.5714
$ws.Range("H33").Value = 6856.5454
$ws.Range("I33").Value = 926.8570999999999
$ws.Range("K33").Value = 926.8570999999999
$ws.Range("M33").Value = -547.8570999999999
$ws.Range("H35").Value = 1750
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 2500
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 2500
$ws.Range("M35").Value = -706
$ws.Range("N35").Value = -3088
$ws.Range("H36").Value = 6500
$ws.Range("I36").Value = 6500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 6500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -6112
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 40000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 40000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 40000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -40782
$ws.Range("H40").Value = 6500
$ws.Range("I40").Value = 6500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6340
$ws.Range("N40").ClearContents()
$ws.Range("H44").Value = 27500
$ws.Range("I44").Value = 25000
$ws.Range("K44").Value = 25000
$ws.Range("M44").Value = -24558
$ws.Range("H49").Value = 40000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 40000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 40000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -40364
$ws.Range("H56").Value = 26666.666
$ws.Range("I56").Value = 26666.666
$ws.Range("K56").Value = 26666.666
$ws.Range("M56").Value = -25821.666
$ws.Range("H60").Value = 27644
$ws.Range("I60").Value = 13092
$ws.Range("J60").Value = 31282
$ws.Range("K60").Value = 13092
$ws.Range("L60").Value = 31282
$ws.Range("M60").Value = -12581
$ws.Range("N60").Value = -32304

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 713.8333
$ws.Range("I5").Value = 489
$ws.Range("K5").Value = 1467
$ws.Range("M5").Value = -1355
$ws.Range("H135").Value = 713.8333
$ws.Range("I135").Value = 489
$ws.Range("K135").Value = 4401
$ws.Range("M135").Value = -1866

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1398.7778
$ws.Range("I97").Value = 1509.75
$ws.Range("K97").Value = 1509.75
$ws.Range("M97").Value = -1013.75
$ws.Range("H122").Value = 4452
$ws.Range("I122").Value = 4058.4
$ws.Range("K122").Value = 12175.2
$ws.Range("M122").Value = -9725.200000000001
$ws.Range("H126").Value = 9850.111000000001
$ws.Range("I126").Value = 5108.5
$ws.Range("K126").Value = 15325.5
$ws.Range("M126").Value = -12855.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2389.5
$ws.Range("I22").Value = 1749.75
$ws.Range("J22").Value = 2816
$ws.Range("K22").Value = 1749.75
$ws.Range("L22").Value = 2816
$ws.Range("M22").Value = -1454.75
$ws.Range("N22").Value = -3406
$ws.Range("H27").Value = 2389.5
$ws.Range("I27").Value = 1749.75
$ws.Range("J27").Value = 2816
$ws.Range("K27").Value = 1749.75
$ws.Range("L27").Value = 2816
$ws.Range("M27").Value = -1642.75
$ws.Range("N27").Value = -3030
$ws.Range("H46").Value = 169191.25
$ws.Range("I46").Value = 400719.6
$ws.Range("J46").Value = 3813.8572
$ws.Range("K46").Value = 400719.6
$ws.Range("L46").Value = 3813.8572
$ws.Range("M46").Value = -400531.6
$ws.Range("N46").Value = -4189.8572
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 2899.6667
$ws.Range("I136").Value = 2899.6667
$ws.Range("K136").Value = 8699.000100000001
$ws.Range("M136").Value = -6149.000100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 966.6667
$ws.Range("I17").Value = 950
$ws.Range("K17").Value = 950
$ws.Range("M17").Value = -778
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H52").Value = 46666.668
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40452
$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 3000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 3000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -25608
$ws.Range("H96").Value = 28566.25
$ws.Range("J96").Value = 100000
$ws.Range("L96").Value = 100000
$ws.Range("N96").Value = -102746
$ws.Range("H132").Value = 2533.3333
$ws.Range("I132").Value = 2533.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7599.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5069.999899999999
$ws.Range("N132").ClearContents()
